# Case and Fatality Demographics Data Updated
# Updates the three "Fatalities by ..." sheets with the latest counts
# (workbook refreshed from 10.14.21 source pull to 10.21.21 pull).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Fatalities by Age Group
# ---------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")

$wsAge.Range("B4").Value  = 70
$wsAge.Range("B5").Value  = 584
$wsAge.Range("B6").Value  = 1885
$wsAge.Range("B7").Value  = 4714
$wsAge.Range("B8").Value  = 9029
$wsAge.Range("B9").Value  = 6939
$wsAge.Range("B10").Value = 8259
$wsAge.Range("B11").Value = 8849
$wsAge.Range("B12").Value = 8413
$wsAge.Range("B13").Value = 20019
$wsAge.Range("B15").Value = 68792

# ---------------------------------------------------------------
# Fatalities by Gender
# ---------------------------------------------------------------
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")

$wsGender.Range("B2").Value = 28753
$wsGender.Range("B3").Value = 40038

# ---------------------------------------------------------------
# Fatalities by Race-Ethnicity
# ---------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")

$wsRace.Range("B2").Value = 1283
$wsRace.Range("B3").ClearFormats()
$wsRace.Range("B3").Value = 7175
$wsRace.Range("B4").Value = 30179
$wsRace.Range("B5").Value = 406
$wsRace.Range("B6").Value = 29709

# Drop the two trailing blank/style-only rows below the table.
$wsRace.Range("A10:O11").Delete()

# ---------------------------------------------------------------
# Restore selections on each touched sheet, finishing on the sheet
# that should remain the active tab.
# ---------------------------------------------------------------
$wsGender.Activate()
$wsGender.Range("B2:B4").Select()

$wsRace.Activate()
$wsRace.Range("D7").Select()

$wsAge.Activate()
$wsAge.Range("G10").Select()
